$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version property value (row 3)
$ws.Range("B3").Value = "0.2.0"

# Update Date property value (row 8)
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# Insert a new "Jurisdiction" row right after the "Contact" row (row 10),
# pushing Description/Purpose/Copyright/Source/Target down by one row.
$ws.Rows.Item(11).EntireRow.Insert()

# Copy the formatting of the row that is now below (the old row 11, now row 12)
# onto the freshly inserted row 11 so it keeps the same borders/alignment.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
